# re-pick test cases; refine plots
#
# 1) Rename the existing sheet "Tabelle1" -> "evaluation"
# 2) Add a new sheet "selection" right after it, with a grid of
#    test-case selection notes / small tally tables
# 3) Update selections on both sheets

$wb = $excel.ActiveWorkbook

# --- sheet1: Tabelle1 -> evaluation -------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "evaluation"

# --- sheet2: new "selection" sheet, inserted right after evaluation -----
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "selection"

# Header / label block (A:E) -------------------------------------------
$ws2.Range("A1").Value = "frequency"
$ws2.Range("B1").Value = "stability"
$ws2.Range("C1").Value = "stable"
$ws2.Range("D1").Value = "ghosting"
$ws2.Range("E1").Value = "lemma"

$ws2.Range("C2").Value = "topical"
$ws2.Range("D2").Value = "poppygate"

$ws2.Range("B3").Value = "trend"
$ws2.Range("C3").Value = "increasing"
$ws2.Range("D3").Value = "ghosting"

$ws2.Range("C4").Value = "decreasing"
$ws2.Range("D4").Value = "robo-signing"
$ws2.Range("E4").Value = "ghosting"

$ws2.Range("B5").Value = "time window"
$ws2.Range("C5").Value = "short"
$ws2.Range("D5").Value = "robo-signing"
$ws2.Range("E5").Value = "poppygate"

$ws2.Range("C6").Value = "long"
$ws2.Range("D6").Value = "ghosting"
$ws2.Range("E6").Value = "robo-signing"

$ws2.Range("E7").Value = "hyperlocal"

$ws2.Range("A8").Value = "centralization"
$ws2.Range("B8").Value = "trend"
$ws2.Range("C8").Value = "decreasing"
$ws2.Range("D8").Value = "hyperlocal"
$ws2.Range("E8").Value = "solopreneur"

$ws2.Range("C9").Value = "increasing"
$ws2.Range("D9").Value = "solopreneur, robo-signing"

$ws2.Range("B10").Value = "overall"
$ws2.Range("C10").Value = "low"
$ws2.Range("D10").Value = "ghosting"
$ws2.Range("E10").Value = "newsjacking"

$ws2.Range("C11").Value = "high"
$ws2.Range("D11").Value = "alt-left"

# Small tally tables (F:P), centered values ------------------------------
$ws2.Range("F1").Value = "freq"
$ws2.Range("L1").Value = "centralization"

$ws2.Range("F2").Value = "stability"
$ws2.Range("H2").Value = "trend"
$ws2.Range("J2").Value = "window"
$ws2.Range("L2").Value = "trend"
$ws2.Range("O2").Value = "overall"

$ws2.Range("F3").Value = "stable"
$ws2.Range("G3").Value = "unstable"
$ws2.Range("H3").Value = "increasing"
$ws2.Range("I3").Value = "decreasing"
$ws2.Range("J3").Value = "short"
$ws2.Range("K3").Value = "long"
$ws2.Range("L3").Value = "increasing"
$ws2.Range("M3").Value = "stable"
$ws2.Range("N3").Value = "decreasing"
$ws2.Range("O3").Value = "low"
$ws2.Range("P3").Value = "high"

$ws2.Range("F4").Value = 2
$ws2.Range("H4").Value = 2
$ws2.Range("K4").Value = 1
$ws2.Range("M4").Value = 1

$ws2.Range("G5").Value = 2

$ws2.Range("I6").Value = 1
$ws2.Range("J6").Value = 1
$ws2.Range("N6").Value = 2

$ws2.Range("I7").Value = 2
$ws2.Range("L7").Value = 1

$ws2.Range("N10").Value = 1

# center-align the tally cells (F:P block), matching the sheet's look
$centeredCells = @( `
  "F1","L1", `
  "F2","H2","J2","L2","O2", `
  "F3","G3","H3","I3","J3","K3","L3","M3","N3","O3","P3", `
  "F4","H4","K4","M4", `
  "G5", `
  "I6","J6","N6", `
  "I7","L7", `
  "N10" `
)
foreach ($cellRef in $centeredCells) {
  $ws2.Range($cellRef).HorizontalAlignment = -4108   # xlCenter
}

try {
  $win = $excel.ActiveWindow
  $win.Zoom = 163
} catch {
}

# --- selections -----------------------------------------------------
$ws2.Range("J11").Select()
try {
  $win = $excel.ActiveWindow
  $win.ScrollColumn = 4
  $win.ScrollRow = 1
} catch {
}

$ws1.Range("E15").Select()
